# Update market-board derived profit figures across multiple sheets
# (values refreshed by scheduled data-scrape runner)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 3575.5
$ws.Range("I100").Value = 3399.2856
$ws.Range("J100").Value = 3986.6667
$ws.Range("K100").Value = 3399.2856
$ws.Range("L100").Value = 3986.6667
$ws.Range("M100").Value = -2858.2856
$ws.Range("N100").Value = -5068.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H124").Value = 40000
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 40000
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 40000
$ws.Range("N124").Value = -49820

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 838.5
$ws.Range("I135").Value = 697.8570999999999
$ws.Range("J135").Value = 1166.6666
$ws.Range("K135").Value = 6280.7139
$ws.Range("L135").Value = 10499.9994
$ws.Range("M135").Value = -3745.7139
$ws.Range("N135").Value = -15569.9994

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 10001996
$ws.Range("I2").Value = 31252238
$ws.Range("J2").Value = 1882.3529
$ws.Range("K2").Value = 31252238
$ws.Range("L2").Value = 1882.3529
$ws.Range("M2").Value = -31252125
$ws.Range("N2").Value = -2108.3529

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6913.3433
$ws.Range("I32").Value = 6126.967
$ws.Range("J32").Value = 14777.111
$ws.Range("K32").Value = 6126.967
$ws.Range("L32").Value = 14777.111
$ws.Range("M32").Value = -5839.967
$ws.Range("N32").Value = -15351.111

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3794.9167
$ws.Range("I61").Value = 4399.25
$ws.Range("J61").Value = 3492.75
$ws.Range("K61").Value = 4399.25
$ws.Range("L61").Value = 3492.75
$ws.Range("M61").Value = -4187.25
$ws.Range("N61").Value = -3916.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1264.175
$ws.Range("I74").Value = 833.6857
$ws.Range("J74").Value = 4277.6
$ws.Range("K74").Value = 833.6857
$ws.Range("L74").Value = 4277.6
$ws.Range("M74").Value = 40.3143
$ws.Range("N74").Value = -6025.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1264.175
$ws.Range("I77").Value = 833.6857
$ws.Range("J77").Value = 4277.6
$ws.Range("K77").Value = 4168.4285
$ws.Range("L77").Value = 21388
$ws.Range("M77").Value = 199.5715
$ws.Range("N77").Value = -30124

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 10001996
$ws.Range("I116").Value = 31252238
$ws.Range("J116").Value = 1882.3529
$ws.Range("K116").Value = 31252238
$ws.Range("L116").Value = 1882.3529
$ws.Range("M116").Value = -31249944
$ws.Range("N116").Value = -6470.3529

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2008.289
$ws.Range("I122").Value = 1654.9688
$ws.Range("J122").Value = 2878
$ws.Range("K122").Value = 4964.9064
$ws.Range("L122").Value = 8634
$ws.Range("M122").Value = -2514.9064
$ws.Range("N122").Value = -13534

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3794.9167
$ws.Range("I136").Value = 4399.25
$ws.Range("J136").Value = 3492.75
$ws.Range("K136").Value = 13197.75
$ws.Range("L136").Value = 10478.25
$ws.Range("M136").Value = -10647.75
$ws.Range("N136").Value = -15578.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 10001996
$ws.Range("I3").Value = 31252238
$ws.Range("J3").Value = 1882.3529
$ws.Range("K3").Value = 31252238
$ws.Range("L3").Value = 1882.3529
$ws.Range("M3").Value = -31252124
$ws.Range("N3").Value = -2110.3529

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3383.9487
$ws.Range("I134").Value = 3754.76
$ws.Range("J134").Value = 2721.7856
$ws.Range("K134").Value = 11264.28
$ws.Range("L134").Value = 8165.3568
$ws.Range("M134").Value = -8729.280000000001
$ws.Range("N134").Value = -13235.3568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 7203.875
$ws.Range("I134").Value = 9174.799999999999
$ws.Range("J134").Value = 3919
$ws.Range("K134").Value = 27524.4
$ws.Range("L134").Value = 11757
$ws.Range("M134").Value = -24989.4
$ws.Range("N134").Value = -16827

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 1994.8889
$ws.Range("I118").Value = 271.6
$ws.Range("J118").Value = 2657.6924
$ws.Range("K118").Value = 814.8000000000001
$ws.Range("L118").Value = 7973.0772
$ws.Range("M118").Value = 428.1999999999999
$ws.Range("N118").Value = -10459.0772

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 16999.7
$ws.Range("I120").Value = 10000
$ws.Range("J120").Value = 19999.572
$ws.Range("K120").Value = 30000
$ws.Range("L120").Value = 59998.716
$ws.Range("M120").Value = -25162
$ws.Range("N120").Value = -69674.716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4728
$ws.Range("I70").Value = 4956
$ws.Range("J70").Value = 4500
$ws.Range("K70").Value = 4956
$ws.Range("L70").Value = 4500
$ws.Range("M70").Value = -4686
$ws.Range("N70").Value = -5040

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4728
$ws.Range("I73").Value = 4956
$ws.Range("J73").Value = 4500
$ws.Range("K73").Value = 4956
$ws.Range("L73").Value = 4500
$ws.Range("M73").Value = -4020
$ws.Range("N73").Value = -6372

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3663.08
$ws.Range("I122").Value = 3957.7
$ws.Range("J122").Value = 3466.6667
$ws.Range("K122").Value = 11873.1
$ws.Range("L122").Value = 10400.0001
$ws.Range("M122").Value = -9423.099999999999
$ws.Range("N122").Value = -15300.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4304.1
$ws.Range("I132").Value = 5540.3076
$ws.Range("J132").Value = 3358.7646
$ws.Range("K132").Value = 16620.9228
$ws.Range("L132").Value = 10076.2938
$ws.Range("M132").Value = -14090.9228
$ws.Range("N132").Value = -15136.2938

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2287.5
$ws.Range("I61").Value = 745
$ws.Range("J61").Value = 10000
$ws.Range("K61").Value = 745
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = -543
$ws.Range("N61").Value = -10404

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2287.5
$ws.Range("I113").Value = 745
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 745
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = 1425
$ws.Range("N113").Value = -14340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8380.348
$ws.Range("I132").Value = 3249.7144
$ws.Range("J132").Value = 10625
$ws.Range("K132").Value = 9749.143199999999
$ws.Range("L132").Value = 31875
$ws.Range("M132").Value = -7219.143199999999
$ws.Range("N132").Value = -36935

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 550
$ws.Range("I107").Value = 550
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1650
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 270
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3562.125
$ws.Range("I122").Value = 2099.6
$ws.Range("J122").Value = 5999.6665
$ws.Range("K122").Value = 6298.799999999999
$ws.Range("L122").Value = 17998.9995
$ws.Range("M122").Value = -3848.799999999999
$ws.Range("N122").Value = -22898.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 17062.486
$ws.Range("I132").Value = 2463.4644
$ws.Range("J132").Value = 75458.57000000001
$ws.Range("K132").Value = 7390.3932
$ws.Range("L132").Value = 226375.71
$ws.Range("M132").Value = -4860.3932
$ws.Range("N132").Value = -231435.71

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2438.2415
$ws.Range("I136").Value = 1914.1818
$ws.Range("J136").Value = 4085.2856
$ws.Range("K136").Value = 5742.5454
$ws.Range("L136").Value = 12255.8568
$ws.Range("M136").Value = -3192.5454
$ws.Range("N136").Value = -17355.8568

